$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 514.7931
$ws.Range("I19").Value = 524.9167
$ws.Range("J19").Value = 507.64706
$ws.Range("K19").Value = 524.9167
$ws.Range("L19").Value = 507.64706
$ws.Range("M19").Value = -349.9167
$ws.Range("N19").Value = -857.64706

$ws.Range("H28").Value = 840.5
$ws.Range("I28").Value = 860.75
$ws.Range("J28").Value = 800
$ws.Range("K28").Value = 860.75
$ws.Range("L28").Value = 800
$ws.Range("M28").Value = -375.75
$ws.Range("N28").Value = -1770

$ws.Range("H49").Value = 1237.5
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 1237.5
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 3712.5
$ws.Range("N49").Value = -3984.5
$ws.Range("M49").ClearContents()

$ws.Range("H51").Value = 4825.0625
$ws.Range("I51").Value = 5980.2
$ws.Range("J51").Value = 4300
$ws.Range("K51").Value = 5980.2
$ws.Range("L51").Value = 4300
$ws.Range("M51").Value = -5496.2
$ws.Range("N51").Value = -5268

$ws.Range("H116").Value = 143959
$ws.Range("I116").Value = 236598.33
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 236598.33
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -233156.33
$ws.Range("N116").Value = -11884

$ws.Range("H125").Value = 566.8421
$ws.Range("I125").Value = 626
$ws.Range("J125").Value = 465.42856
$ws.Range("K125").Value = 5634
$ws.Range("L125").Value = 4188.85704
$ws.Range("M125").Value = -3174
$ws.Range("N125").Value = -9108.857039999999

$ws.Range("H137").Value = 2734.162
$ws.Range("I137").Value = 2852.1155
$ws.Range("J137").Value = 2455.3635
$ws.Range("K137").Value = 8556.3465
$ws.Range("L137").Value = 7366.0905
$ws.Range("M137").Value = -6006.3465
$ws.Range("N137").Value = -12466.0905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3578.99
$ws.Range("I32").Value = 3062.7659
$ws.Range("J32").Value = 11666.5
$ws.Range("K32").Value = 3062.7659
$ws.Range("L32").Value = 11666.5
$ws.Range("M32").Value = -2775.7659
$ws.Range("N32").Value = -12240.5

$ws.Range("H61").Value = 3945.0322
$ws.Range("I61").Value = 2723.2942
$ws.Range("J61").Value = 5428.5713
$ws.Range("K61").Value = 2723.2942
$ws.Range("L61").Value = 5428.5713
$ws.Range("M61").Value = -2511.2942
$ws.Range("N61").Value = -5852.5713

$ws.Range("H88").Value = 2640
$ws.Range("I88").Value = 2633.3333
$ws.Range("J88").Value = 2666.6667
$ws.Range("K88").Value = 2633.3333
$ws.Range("L88").Value = 2666.6667
$ws.Range("M88").Value = -2227.3333
$ws.Range("N88").Value = -3478.6667

$ws.Range("H91").Value = 2640
$ws.Range("I91").Value = 2633.3333
$ws.Range("J91").Value = 2666.6667
$ws.Range("K91").Value = 2633.3333
$ws.Range("L91").Value = 2666.6667
$ws.Range("M91").Value = -1229.3333
$ws.Range("N91").Value = -5474.6667

$ws.Range("H97").Value = 1989.5
$ws.Range("I97").Value = 2067.5
$ws.Range("J97").Value = 1755.5
$ws.Range("K97").Value = 2067.5
$ws.Range("L97").Value = 1755.5
$ws.Range("M97").Value = -1571.5
$ws.Range("N97").Value = -2747.5

$ws.Range("H113").Value = 35716
$ws.Range("J113").Value = 35716
$ws.Range("L113").Value = 35716
$ws.Range("N113").Value = -44394

$ws.Range("H132").Value = 2651.1333
$ws.Range("I132").Value = 2405.9534
$ws.Range("J132").Value = 3271.2942
$ws.Range("K132").Value = 7217.860199999999
$ws.Range("L132").Value = 9813.882599999999
$ws.Range("M132").Value = -4687.860199999999
$ws.Range("N132").Value = -14873.8826

$ws.Range("H136").Value = 3945.0322
$ws.Range("I136").Value = 2723.2942
$ws.Range("J136").Value = 5428.5713
$ws.Range("K136").Value = 8169.882599999999
$ws.Range("L136").Value = 16285.7139
$ws.Range("M136").Value = -5619.882599999999
$ws.Range("N136").Value = -21385.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H107").Value = 1918.5555
$ws.Range("I107").Value = 2061.4
$ws.Range("J107").Value = 1204.3334
$ws.Range("K107").Value = 2061.4
$ws.Range("L107").Value = 1204.3334
$ws.Range("M107").Value = -141.4000000000001
$ws.Range("N107").Value = -5044.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3381.7017
$ws.Range("I31").Value = 2820.28
$ws.Range("J31").Value = 3820.3125
$ws.Range("K31").Value = 2820.28
$ws.Range("L31").Value = 3820.3125
$ws.Range("M31").Value = -2525.28
$ws.Range("N31").Value = -4410.3125

$ws.Range("H34").Value = 3381.7017
$ws.Range("I34").Value = 2820.28
$ws.Range("J34").Value = 3820.3125
$ws.Range("K34").Value = 2820.28
$ws.Range("L34").Value = 3820.3125
$ws.Range("M34").Value = -2618.28
$ws.Range("N34").Value = -4224.3125

$ws.Range("H122").Value = 3162.4
$ws.Range("I122").Value = 3453
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 10359
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -7909
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 1762.0986
$ws.Range("I132").Value = 1051.94
$ws.Range("J132").Value = 3452.9524
$ws.Range("K132").Value = 3155.82
$ws.Range("L132").Value = 10358.8572
$ws.Range("M132").Value = -625.8200000000002
$ws.Range("N132").Value = -15418.8572

$ws.Range("H141").Value = 30065.555
$ws.Range("J141").Value = 28823.75
$ws.Range("L141").Value = 28823.75
$ws.Range("N141").Value = -39183.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 2901.3635
$ws.Range("I121").Value = 5685
$ws.Range("J121").Value = 2461.842
$ws.Range("K121").Value = 17055
$ws.Range("L121").Value = 7385.526
$ws.Range("M121").Value = -15745
$ws.Range("N121").Value = -10005.526

$ws.Range("H131").Value = 1295.4648
$ws.Range("I131").Value = 647.2105
$ws.Range("J131").Value = 1532.3269
$ws.Range("K131").Value = 1941.6315
$ws.Range("L131").Value = 4596.9807
$ws.Range("M131").Value = 3098.3685
$ws.Range("N131").Value = -14676.9807

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 894.75
$ws.Range("I113").Value = 772.7143
$ws.Range("J113").Value = 1749
$ws.Range("K113").Value = 772.7143
$ws.Range("L113").Value = 1749
$ws.Range("M113").Value = 1397.2857
$ws.Range("N113").Value = -6089

$ws.Range("H132").Value = 2961.9106
$ws.Range("I132").Value = 2947
$ws.Range("J132").Value = 2990.9473
$ws.Range("K132").Value = 8841
$ws.Range("L132").Value = 8972.841899999999
$ws.Range("M132").Value = -6311
$ws.Range("N132").Value = -14032.8419

$ws.Range("H133").Value = 29926.666
$ws.Range("J133").Value = 29926.666
$ws.Range("L133").Value = 29926.666
$ws.Range("N133").Value = -40046.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 426
$ws.Range("I113").Value = 384.66666
$ws.Range("J113").Value = 550
$ws.Range("K113").Value = 1153.99998
$ws.Range("L113").Value = 1650
$ws.Range("M113").Value = 1016.00002
$ws.Range("N113").Value = -5990

$ws.Range("H132").Value = 1806.0322
$ws.Range("I132").Value = 986.1818
$ws.Range("J132").Value = 3810.111
$ws.Range("K132").Value = 2958.5454
$ws.Range("L132").Value = 11430.333
$ws.Range("M132").Value = -428.5454
$ws.Range("N132").Value = -16490.333

$ws.Range("H136").Value = 1873.2543
$ws.Range("I136").Value = 1057.5714
$ws.Range("K136").Value = 3172.7142
$ws.Range("M136").Value = -622.7142000000003
